# RTM / Traceability-Matrix.xlsx update
# Adds the "Client" requirement/test-case traceability rows (SRS_Client_001..015
# mapped to TC_Client_001..018) into rows 36-50 of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A: SRS IDs, column B: matching Test Case IDs (some cells combine
# several TC ids separated by long runs of spaces, mirroring the source data).
$data = @(
    @{ Row = 36; A = "SRS_Client_001";  B = "TC_Client_001" },
    @{ Row = 37; A = "SRS_Client_002";  B = "TC_Client_002" },
    @{ Row = 38; A = "SRS_Client_003";  B = "TC_Client_001" },
    @{ Row = 39; A = "SRS_Client_004";  B = "TC_Client_003" },
    @{ Row = 40; A = "SRS_Client_005";  B = "TC_Client_004                            TC_Client_005                             TC_Client_006" },
    @{ Row = 41; A = "SRS_Client_006";  B = "TC_Client_007" },
    @{ Row = 42; A = "SRS_Client_007";  B = "TC_Client_008" },
    @{ Row = 43; A = "SRS_Client_008";  B = "TC_Client_010                               TC_Client_011                       TC_Client_012" },
    @{ Row = 44; A = "SRS_Client_009";  B = "TC_Client_013" },
    @{ Row = 45; A = "SRS_Client_0010"; B = "TC_Client_014" },
    @{ Row = 46; A = "SRS_Client_011";  B = "TC_Client_009" },
    @{ Row = 47; A = "SRS_Client_012";  B = "TC_Client_015                             TC_Client_016" },
    @{ Row = 48; A = "SRS_Client_013";  B = "TC_Client_017" },
    @{ Row = 49; A = "SRS_Client_014";  B = "TC_Client_018" },
    @{ Row = 50; A = "SRS_Client_015";  B = "TC_Client_018" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
}

# Wrapped multi-line-looking cells (those combining several TC ids) use the
# same word-wrap style already present on the sheet (style index 7 => the
# font3/border1/wrapText style), matching B40, B43 and B47.
$ws.Range("B40").WrapText = $true
$ws.Range("B43").WrapText = $true
$ws.Range("B47").WrapText = $true

# Leave the selection where the author left off after data entry.
$ws.Range("A51").Select()

$wb.Save()
